$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("EnergyTransport")
$c = $ws2.Range("N46").Comment
try {
  $c.Shape.TopLeftCell = $ws2.Range("N50")
  Write-Host "ok move"
} catch {
  Write-Host "err: $_"
}
Write-Host "author: $($c.Author())"
